$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Fecha Actualización" date for the row whose "Situación" is "Trabajando"
$ws.Range("F35").Value = 44172

# Filter the table so only rows with "Situación" = "Trabajando" are shown
$tbl = $ws.ListObjects.Item("Tabla1")
$tbl.Range.AutoFilter(4, @("Trabajando"), 7) | Out-Null

# Update the view: scroll the frozen pane so row 8 is the first visible row
# below the header, and select cell F51
$ws.Range("F51").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 8
